$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.709.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.128.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.121.47"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("E10").Value = "  +14.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("E13").Value = "  +5.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.63%  "
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.646.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.608.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.128.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.732"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("E24").Value = "  -3.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.02%  "
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.28%  "
$ws.Range("E33").Value = "  -2.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0877"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.01%  "
$ws.Range("E35").Value = "  +8.50%  "
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +16.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "51.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "454.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.903.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.279"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("E45").Value = "  +1.94%  "
$ws.Range("E46").Value = "  +3.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "128.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.88%  "
